# Working Checkpoint. Reading String from Template and Replacing Chars for Body
#
# Recreates the recruiter contact list: fills in company / contact / email /
# position data, adds mailto hyperlinks for the e-mail column, and relabels
# the header row (First Name / Last Name / Company / Mail ID / Position
# Applying For).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Company column (C) -----------------------------------------------
$ws.Range("C2").Value = "NC State"
$ws.Range("C3").Value = "Google"
$ws.Range("C4").Value = "Microsoft"

# --- First name column (A) ---------------------------------------------
$ws.Range("A2").Value = "abc"
$ws.Range("A3").Value = "def"
$ws.Range("A4").Value = "ghi"

# --- Mail ID column (D), as hyperlinks ----------------------------------
$ws.Range("D2").Value = "abc@gmail.com"
$ws.Range("D3").Value = "def@gmail.com"
$ws.Range("D4").Value = "ghi@gmail.com"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:abc@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:def@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:ghi@gmail.com") | Out-Null

# --- Position Applying For column (E) header ----------------------------
$ws.Range("E1").Value = "Position Applying For"

# --- Header row relabeling (A1/B1) --------------------------------------
$ws.Range("A1").Value = "First Name"
$ws.Range("B1").Value = "Last Name"

# --- Position Applying For column (E) data ------------------------------
$ws.Range("E2").Value = "abc Engineer"
$ws.Range("E3").Value = "def Engineer"
$ws.Range("E4").Value = "ghi Engineer"

# --- Header row relabeling (C1/D1), reusing existing shared strings -----
$ws.Range("C1").Value = "Company"
$ws.Range("D1").Value = "Mail ID"

# --- Column widths for the new Mail ID / Position Applying For columns --
$ws.Columns.Item(4).ColumnWidth = 22.83
$ws.Columns.Item(5).ColumnWidth = 17.5

# --- Final selection, matching the saved state of the workbook ---------
$ws.Range("B3").Select() | Out-Null

Write-Host "Recruiters sheet populated"
